$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "canonical SMILES"
$ws.Range("D3").Value = "c1ccc2c(c1)c([nH+]cn2)Nc3cccc(c3)Cl"
$ws.Range("D4").Value = "c1ccc2c(c1)c(=Nc3cccc(c3)Cl)[nH]cn2"
$ws.Range("D5").Value = "c1ccc2c(c1)c(nc[nH+]2)Nc3cccc(c3)Cl"
$ws.Range("D6").Value = "c1ccc2c(c1)c(ncn2)[NH2+]c3cccc(c3)Cl"
$ws.Range("D7").Value = "c1ccc2c(c1)c(ncn2)[N-]c3cccc(c3)Cl"
$ws.Range("D8").Value = "c1ccc2c(c1)c(nc[nH+]2)[NH2+]c3cccc(c3)Cl"
$ws.Range("D9").Value = "c1ccc2c(c1)c(=Nc3cccc(c3)Cl)nc[nH]2"
$ws.Range("D10").Value = "c1ccc2c(c1)c(ncn2)Nc3cccc(c3)Cl"
$ws.Range("D11").Value = "c1ccc2c(c1)c([nH+]c[nH+]2)Nc3cccc(c3)Cl"
$ws.Range("D12").Value = "c1ccc2c(c1)c([nH+]cn2)[NH2+]c3cccc(c3)Cl"
$ws.Range("D13").Value = "c1ccc2c(c1)c([nH+]c[nH+]2)[NH2+]c3cccc(c3)Cl"

$ws.Columns.Item(4).ColumnWidth = 36
